# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 85 on the active sheet, pushing
# all the existing rows (old 85..167) down by one (new 86..168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 85 - everything below shifts down.
$ws.Rows(85).Insert()

# Populate the newly inserted row 85 with the new weekly record.
$ws.Range("A85").Value = 7
$ws.Range("B85").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C85").Value = "Ñuble"
$ws.Range("D85").Value = 44546
$ws.Range("E85").Value = 16
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100104
$ws.Range("H85").Value = "Frutos de pepita"
$ws.Range("I85").Value = 100104005
$ws.Range("J85").Value = "Pera"
$ws.Range("K85").Value = "Packham's Triumph"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 100
$ws.Range("N85").Value = 10000
$ws.Range("O85").Value = 11000
$ws.Range("P85").Value = 10500
$ws.Range("Q85").Value = '$/caja 16 kilos empedrada'
$ws.Range("R85").Value = "Provincia de Curicó"
$ws.Range("S85").Value = 656
$ws.Range("T85").Value = 16
